# "Updated presentation and demos"
#
# 1) The little "cells react to each other" demo: A2 changes from 6 to 2
#    (B2 stays 3, and the C2 = A2+B2 formula recalculates from 9 to 5).
# 2) The "Observable" / "Observer or Subscriber" demo row (row 3) is
#    removed entirely - those two labels are no longer used anywhere in
#    the workbook, so clearing the row also drops the now-unused shared
#    strings on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the reactive value demo.
$ws.Range("A2").Value = 2

# Remove the now-unused Observable/Observer row without shifting the
# rows below it.
$ws.Range("A3:C3").ClearContents()
